# Bugfix: rename "HL" / "IR" sheets to "EL" / "IE" (country code fixes
# Greece / Ireland) and update generation data for files that now pull
# from separate per-country databases.

$wb = $excel.ActiveWorkbook

# --- Rename mis-labelled sheets ---------------------------------------
$wb.Worksheets.Item("HL").Name = "EL"
$wb.Worksheets.Item("IR").Name = "IE"

# --- Cell value corrections --------------------------------------------
function Set-Val($sheetName, $cellRef, $value) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range($cellRef).Value = $value
}

# AT
Set-Val "AT" "G11" 0.08369999999999989
Set-Val "AT" "G12" 0.4649999999999999
Set-Val "AT" "G13" 0.8319999999999999

# DE
Set-Val "DE" "E4"  3.056189999999999
Set-Val "DE" "G4"  0.002900000000000125
Set-Val "DE" "E6"  2.428600000000001
Set-Val "DE" "G6"  1.571700000000001
Set-Val "DE" "E7"  8.380099999999997
Set-Val "DE" "L7"  0.3613
Set-Val "DE" "L8"  0.6460999999999988
Set-Val "DE" "E9"  1.687400000000004
Set-Val "DE" "G9"  0.9992099999999988
Set-Val "DE" "J9"  0.3774790000000006
Set-Val "DE" "E10" 0.2593999999999994
Set-Val "DE" "G10" 4.118400000000002
Set-Val "DE" "I10" 0.1108
Set-Val "DE" "J10" 0.5912199999999999
Set-Val "DE" "G11" 1.753159000000002
Set-Val "DE" "I11" 0.07519999999999971
Set-Val "DE" "J11" 0.4382499999999983
Set-Val "DE" "D12" 0.26985
Set-Val "DE" "E12" 0.05059999999999931
Set-Val "DE" "G12" 2.987100000000002
Set-Val "DE" "D13" 0.01509999999999989
Set-Val "DE" "E13" 7.831700000000001
Set-Val "DE" "G13" 2.482799999999997
Set-Val "DE" "I13" 0.15896
Set-Val "DE" "J13" 0.4084500000000002
Set-Val "DE" "G14" 1.567499999999999
Set-Val "DE" "J14" 1.294461

# EL (formerly HL)
Set-Val "EL" "G13" 0.8649999999999993

# IT
Set-Val "IT" "K7"  0.134
Set-Val "IT" "K8"  0.15
Set-Val "IT" "K9"  0.08900000000000002
Set-Val "IT" "G10" 3.798000000000002
Set-Val "IT" "K10" 0.352
Set-Val "IT" "K11" 0.4480000000000001
Set-Val "IT" "G12" 12.002
Set-Val "IT" "K12" 0.05399999999999983
Set-Val "IT" "G13" 3.568999999999996

# NL
Set-Val "NL" "G13" 2.671999999999999

# BE
Set-Val "BE" "G10" 0
Set-Val "BE" "G11" 0.3763999999999998
Set-Val "BE" "G13" 0.8247999999999998

# PL
Set-Val "PL" "D13" 0.205

# PT
Set-Val "PT" "G12" 0.8260000000000001

# RO
Set-Val "RO" "G13" 0.8850000000000002

# ES
Set-Val "ES" "G11" 11.978166
Set-Val "ES" "G12" 18.88164
Set-Val "ES" "G13" 4.333049999999993

# DK
Set-Val "DK" "E9" 0.7869999999999999

# FR
Set-Val "FR" "E7" 1.785
